$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing "..._Filtered" entry (row 3) to "..._20Filtered"
$ws.Range("A3").Value = "FSR_N1_Stability(5.00lbf)_20Filtered"

# New row 4 label, new row 5 label (written first so the shared-string
# table fills in the same order the source workbook used)
$ws.Range("A4").Value = "FSR_N1_Stability(5.00lbf)_100Filtered"
$ws.Range("A5").Value = "FSR_N1_Stability(5.00lbf)_50Filtered"

# New row 4: 100%-filtered stability data
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 28.24
$ws.Range("D4").Value = 95.02
$ws.Range("E4").Value = 3.18
$ws.Range("F4").Value = 6.8650000000000002
$ws.Range("G4").Value = 379
$ws.Range("H4").Value = "Removed all Error < 100%"

# New row 5: 50%-filtered stability data (partial row)
$ws.Range("B5").Value = 5
$ws.Range("H5").Value = "Removed all Error < 50%"

# New empty styled row 8 (copy formatting only from the header row)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A8:H8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Column A is now wider to fit the longer labels (closest reachable
# width to the source's recomputed best-fit width of 35.85546875)
$ws.Columns("A").ColumnWidth = 35

# Move / update the active selection to H8
$ws.Range("H8").Select() | Out-Null
